$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two task descriptions in column B ("Actualizacion del plan general")
$ws.Range("B2").Value = "Ver video tutorial de GitHub."
$ws.Range("B5").Value = "Instalar y configurar ruby on rails."

# Nudge the stored column widths to match the refreshed plan layout
$ws.Columns.Item(1).ColumnWidth = 11.033333333333365
$ws.Columns.Item(2).ColumnWidth = 39.23333333333336
$ws.Columns.Item(3).ColumnWidth = 1.8529411764705865
$ws.Columns.Item(4).ColumnWidth = 11.033333333333365
$ws.Columns.Item(5).ColumnWidth = 1.8529411764705865
$ws.Columns.Item(6).ColumnWidth = 19.19803921568627

# Move the active selection as recorded in the saved view state
$ws.Range("B6").Select()
